$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testreg4")

# Column C (altLabel) needs to be widened to fit the new multi-line alt-label text.
$ws.Columns.Item(3).ColumnWidth = 25.3

# New row 8 data: id, label, altLabel (multi-line, wrapped), description, notation
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "multi altLabel"
$ws.Range("C8").WrapText = $true
$ws.Range("C8").Value = "altLabel-multi-line1`naltLabel-multi-line2"
$ws.Range("D8").Value = "concept with multi alt-labels"
$ws.Range("E8").Value = "multi"

$ws.Rows.Item(8).RowHeight = 28.8

$ws.Range("D8").Select()
